$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (also updates defined name reference automatically)
$ws.Name = "20161004 FXI"

# Update values in column B
$ws.Range("B2").Value = "20161004 +FXI-161021P37.00"
$ws.Range("B3").Value = "E:\\Datos\\bolsa\\cuenta personal\\analisis de valores\\Trades activos\\Scanning\\20161004"
$ws.Range("B5").Value = 37
$ws.Range("B9").Value = 0.26
$ws.Range("B10").Value = 38.17
$ws.Range("B15").Value = 36
$ws.Range("B16").Value = 0
$ws.Range("B17").Value = "FXI"
$ws.Range("B18").Value = 0.2162
$ws.Range("B19").Value = "20161004 +FXI-161021P38.00"
$ws.Range("B20").Value = "E:\\Datos\\bolsa\\cuenta personal\\analisis de valores\\Trades activos\\Scanning\\20161004"
$ws.Range("B21").Value = 38
$ws.Range("B22").Value = -0.54
$ws.Range("B23").Value = 38.215
$ws.Range("B28").Value = 40
$ws.Range("B29").Value = 42
$ws.Range("B30").Value = 0.2162

# Update selection
$ws.Range("B6").Select()
